$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.038.64"
$ws.Range("E2").Value = "  +0.03%  "

Set-TextValue "D3" "3.255.25"
$ws.Range("E3").Value = "  +2.52%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue "D5" "607.56"
$ws.Range("E5").Value = "  +0.65%  "

Set-TextValue "D6" "157.15"
$ws.Range("E6").Value = "  +2.08%  "

$ws.Range("E7").Value = "  +0.08%  "

Set-TextValue "D8" "3.253.38"
$ws.Range("E8").Value = "  +2.44%  "

$ws.Range("E9").Value = "  -0.55%  "

Set-TextValue "D10" "0.160"
$ws.Range("E10").Value = "  +1.28%  "

Set-TextValue "D11" "5.70"
$ws.Range("E11").Value = "  +1.27%  "

Set-TextValue "D12" "0.494"
$ws.Range("E12").Value = "  -2.40%  "

Set-TextValue "D13" "0.0000269"
$ws.Range("E13").Value = "  +1.90%  "

Set-TextValue "D14" "38.50"
$ws.Range("E14").Value = "  +0.78%  "

Set-TextValue "D15" "3.790.85"
$ws.Range("E15").Value = "  +2.57%  "

Set-TextValue "D16" "66.122.61"
$ws.Range("E16").Value = "  +0.10%  "

Set-TextValue "D17" "3.255.66"
$ws.Range("E17").Value = "  +2.46%  "

Set-TextValue "D18" "7.30"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("E19").Value = "  +1.28%  "

Set-TextValue "D20" "498.43"
$ws.Range("E20").Value = "  -1.58%  "

Set-TextValue "D21" "15.27"
$ws.Range("E21").Value = "  +0.16%  "

Set-TextValue "D22" "0.747"
$ws.Range("E22").Value = "  +2.89%  "

Set-TextValue "D23" "8.04"
$ws.Range("E23").Value = "  +0.63%  "

Set-TextValue "D24" "14.58"
$ws.Range("E24").Value = "  -1.11%  "

Set-TextValue "D25" "86.68"
$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +1.31%  "

Set-TextValue "D28" "9.10"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D29" "2.35"
$ws.Range("E29").Value = "  -0.81%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.130"
$ws.Range("E30").Value = "  +44.40%  "

Set-TextValue "D31" "7.02"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  -6.35%  "

Set-TextValue "D33" "27.79"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("E35").Value = "  -3.00%  "

Set-TextValue "D36" "6.39"
$ws.Range("E36").Value = "  -1.10%  "

Set-TextValue "D37" "3.46"
$ws.Range("E37").Value = "  +20.74%  "

Set-TextValue "D38" "55.56"
$ws.Range("E38").Value = "  +0.55%  "

Set-TextValue "D39" "0.0₃0785"
$ws.Range("E39").Value = "  +10.80%  "

Set-TextValue "D40" "494.63"
$ws.Range("E40").Value = "  -2.88%  "

Set-TextValue "D41" "0.0421"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("E42").Value = "  +1.67%  "

Set-TextValue "D43" "8.76"
$ws.Range("E43").Value = "  +0.15%  "

Set-TextValue "D44" "2.56"
$ws.Range("E44").Value = "  +4.01%  "

Set-TextValue "D45" "2.995.20"
$ws.Range("E45").Value = "  +5.95%  "

Set-TextValue "D46" "0.290"
$ws.Range("E46").Value = "  -2.46%  "

Set-TextValue "D47" "28.67"
$ws.Range("E47").Value = "  +3.03%  "

$ws.Range("E48").Value = "  +4.13%  "

$ws.Range("E49").Value = "  +2.17%  "

$ws.Range("E50").Value = "  -0.01%  "

Set-TextValue "D51" "121.16"
$ws.Range("E51").Value = "  -2.12%  "
